$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
foreach ($f in $ftr.Range.Fields) {
    $f.Update()
    Write-Host "Field result after update:" $f.Result.Text
}
Write-Host "ComputeStatistics pages:" $d.ComputeStatistics(2)
